$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old import had a bug: data was shifted one column to the right (header
# started at B1, a spurious 0-based index sat in column A), game scores were
# stored as text, and an extra (incorrect) 6th game row was present.
#
# Step 1: while the header format is still intact on B1:F1, extend it onto
# A1 too (re-using the workbook's existing bold/bordered/centered style
# instead of minting a new one).
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Step 2: wipe everything else - all the old data rows (2-6) plus the old F
# column header cell, which won't exist in the corrected 5-column layout.
$ws.Range("A2:F6").Clear()
$ws.Range("F1").Clear()

# ---- Header row (A1:E1) ---------------------------------------------------
$ws.Range("A1").Value = "Date"
$ws.Range("B1").Value = "Home"
$ws.Range("C1").Value = "HomeScore"
$ws.Range("D1").Value = "Visitor"
$ws.Range("E1").Value = "VisitorScore"

# ---- Data rows (A2:E5) -----------------------------------------------------
# Column A holds the game date as text - pre-format it as Text so Excel
# doesn't auto-convert the date-like strings into date serial numbers.
$ws.Range("A2:A5").NumberFormat = "@"

$ws.Cells.Item(2,1).Value = "11/02/1872"
$ws.Cells.Item(2,2).Value = "Rutgers"
$ws.Cells.Item(2,3).Value = 0
$ws.Cells.Item(2,4).Value = "Columbia"
$ws.Cells.Item(2,5).Value = 0

$ws.Cells.Item(3,1).Value = "11/09/1872"
$ws.Cells.Item(3,2).Value = "Columbia"
$ws.Cells.Item(3,3).Value = 5
$ws.Cells.Item(3,4).Value = "Rutgers"
$ws.Cells.Item(3,5).Value = 7

$ws.Cells.Item(4,1).Value = "11/16/1872"
$ws.Cells.Item(4,2).Value = "Columbia"
$ws.Cells.Item(4,3).Value = 0
$ws.Cells.Item(4,4).Value = "Yale"
$ws.Cells.Item(4,5).Value = 3

$ws.Cells.Item(5,1).Value = "11/16/1872"
$ws.Cells.Item(5,2).Value = "Rutgers"
$ws.Cells.Item(5,3).Value = 1
$ws.Cells.Item(5,4).Value = "Princeton"
$ws.Cells.Item(5,5).Value = 4

$wb.Save()
